# "Generate Report for Handoff" — refresh the localization-status report
# with a new handoff batch (new source file GUIDs, reset handoff/handback
# bookkeeping columns, updated timestamps and status).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: set a cell's text value while forcing the *text* type even if
# the string looks like a boolean ("True"/"False") so Excel doesn't
# auto-coerce it into a real Boolean cell. A leading apostrophe forces
# text entry; we then strip the quote-prefix formatting it leaves
# behind so the cell's style matches a plain, unformatted text cell.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    if ($text -eq "True" -or $text -eq "False" -or $text -eq "") {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# Helper: remove any hyperlink anchored at a specific cell address
# (e.g. "$I$2"). Re-scans the live collection after every delete since
# previously-captured hyperlink references go stale once any one of
# them is removed.
# ---------------------------------------------------------------------
function Remove-HyperlinkAt($ws, [string]$addr) {
    $again = $true
    while ($again) {
        $again = $false
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $addr) {
                $hl.Delete()
                $again = $true
                break
            }
        }
    }
}

# ---------------------------------------------------------------------
# Helper: update the display text of the hyperlink anchored at a given
# cell address, leaving its target URL (r:id) untouched.
# ---------------------------------------------------------------------
function Set-HyperlinkDisplay($ws, [string]$addr, [string]$display) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $display
        }
    }
}

# =======================================================================
# Overview sheet
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "b19eeaea-0863-44cf-b821-197297319948.md"
$wsOverview.Range("B2").Value = "e2e\b19eeaea-0863-44cf-b821-197297319948.md"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 15:04:38"

$wsOverview.Range("A3").Value = "ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
$wsOverview.Range("B3").Value = "e2e\ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 15:04:38"

Set-HyperlinkDisplay $wsOverview '$B$2' "e2e\b19eeaea-0863-44cf-b821-197297319948.md"
Set-HyperlinkDisplay $wsOverview '$B$3' "e2e\ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# =======================================================================
# zh-cn sheet
# =======================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "b19eeaea-0863-44cf-b821-197297319948.md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("G2").Value = "b19eeaea-0863-44cf-b821-197297319948.9b2305b8248363bf7ae7e0811b99ddb3ff168fa8.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 15:04:32"
Set-TextValue $wsZhCn.Range("I2") ""
Set-TextValue $wsZhCn.Range("J2") ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Range("A3").Value = "ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsZhCn.Range("F3") "True"
$wsZhCn.Range("G3").Value = "b19eeaea-0863-44cf-b821-197297319948.9b2305b8248363bf7ae7e0811b99ddb3ff168fa8.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 15:04:32"
Set-TextValue $wsZhCn.Range("I3") ""
Set-TextValue $wsZhCn.Range("J3") ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"

Set-HyperlinkDisplay $wsZhCn '$A$2' "b19eeaea-0863-44cf-b821-197297319948.md"
Set-HyperlinkDisplay $wsZhCn '$A$3' "ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
Remove-HyperlinkAt $wsZhCn '$I$2'
Remove-HyperlinkAt $wsZhCn '$I$3'

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# =======================================================================
# de-de sheet
# =======================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "b19eeaea-0863-44cf-b821-197297319948.md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("G2").Value = "b19eeaea-0863-44cf-b821-197297319948.9b2305b8248363bf7ae7e0811b99ddb3ff168fa8.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-16 15:04:38"
Set-TextValue $wsDeDe.Range("I2") ""
Set-TextValue $wsDeDe.Range("J2") ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("A3").Value = "ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsDeDe.Range("F3") "True"
$wsDeDe.Range("G3").Value = "b19eeaea-0863-44cf-b821-197297319948.9b2305b8248363bf7ae7e0811b99ddb3ff168fa8.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 15:04:38"
Set-TextValue $wsDeDe.Range("I3") ""
Set-TextValue $wsDeDe.Range("J3") ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"

Set-HyperlinkDisplay $wsDeDe '$A$2' "b19eeaea-0863-44cf-b821-197297319948.md"
Set-HyperlinkDisplay $wsDeDe '$A$3' "ffffbd438b67-aaf4-447e-b6cb-73263a696983.md"
Remove-HyperlinkAt $wsDeDe '$I$2'
Remove-HyperlinkAt $wsDeDe '$I$3'

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
